# Updates cryptos list figures (prices / 1h volume %) per the Dec 25 2023
# GitHub Actions refresh, including the Cronos/Aave row swap at rows 49-50.
#
# Values that look like plain numbers (e.g. "120.84") get a leading literal
# apostrophe so Excel's COM layer stores them as TEXT (matching the source
# workbook, which keeps every Price/Volume cell as a string) instead of
# silently converting them to numeric cells. The apostrophe itself is a
# formatting marker, not part of the stored value: $ws.Range("D5").Value
# ends up equal to "120.84", not "'120.84".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value
$updates = [ordered]@{
    "D2" = '43.639.37'
    "E2" = '  +0.97%  '
    "D3" = '2.277.16'
    "E3" = '  +0.36%  '
    "E4" = '  -0.12%  '
    "D5" = '''120.84'
    "E5" = '  +6.52%  '
    "D6" = '''267.08'
    "E6" = '  +0.89%  '
    "E7" = '  +5.06%  '
    "E8" = '  +0.11%  '
    "E9" = '  +5.20%  '
    "D10" = '''48.45'
    "E10" = '  +0.88%  '
    "D11" = '''0.0949'
    "E11" = '  +2.61%  '
    "D12" = '''9.33'
    "E12" = '  +6.72%  '
    "E13" = '  -0.66%  '
    "D14" = '''15.69'
    "E14" = '  +2.00%  '
    "D15" = '''0.916'
    "E15" = '  +6.76%  '
    "D16" = '2.621.75'
    "E16" = '  +0.54%  '
    "D17" = '2.276.95'
    "E17" = '  +0.32%  '
    "D18" = '43.658.07'
    "E18" = '  +1.28%  '
    "E19" = '  +3.45%  '
    "D20" = '''6.95'
    "E20" = '  -0.47%  '
    "D21" = '''72.37'
    "E21" = '  +1.75%  '
    "D22" = '''2.41'
    "E22" = '  -0.09%  '
    "D23" = '''235.98'
    "E23" = '  +2.51%  '
    "E24" = '  -3.02%  '
    "D25" = '''2.90'
    "E25" = '  +2.55%  '
    "D26" = '''11.97'
    "E26" = '  +5.96%  '
    "E27" = '  +1.72%  '
    "D28" = '''43.36'
    "E28" = '  +5.40%  '
    "D29" = '''3.41'
    "E29" = '  +1.10%  '
    "E30" = '  +0.60%  '
    "D31" = '''173.57'
    "E31" = '  +1.17%  '
    "D32" = '''21.71'
    "E32" = '  +1.99%  '
    "D33" = '''0.0926'
    "E33" = '  +2.49%  '
    "D34" = '''5.80'
    "E34" = '  +3.27%  '
    "E35" = '  +4.17%  '
    "D36" = '''4.28'
    "E36" = '  +12.01%  '
    "D37" = '''0.0386'
    "E37" = '  +10.29%  '
    "E38" = '  +0.56%  '
    "E39" = '  +5.17%  '
    "E40" = '  +5.42%  '
    "D41" = '''74.07'
    "E41" = '  -1.18%  '
    "E42" = '  -3.28%  '
    "E43" = '  +2.55%  '
    "D44" = '''0.999'
    "E44" = '  -0.33%  '
    "D45" = '''5.94'
    "E45" = '  -2.73%  '
    "E46" = '  +1.79%  '
    "D47" = '''75.86'
    "E47" = '  +45.40%  '
    "E48" = '  +3.88%  '
    "B49" = 'Aave'
    "C49" = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    "D49" = '''103.24'
    "E49" = '  +2.77%  '
    "B50" = 'Cronos'
    "C50" = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    "D50" = '''0.101'
    "E50" = '  +1.70%  '
    "D51" = '''8.49'
    "E51" = '  -1.09%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
